$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ten_lists")

# The list id for this workbook moved from S000 to S018
$ws.Range("F1").Value = "S018"

# Each of the five paired-list blocks (rows 3, 10, 17, 24, 31) gets two new
# "starting" columns (E and K) describing where the DD/SD numbering for that
# block begins, alongside the existing walk/no-walk category labels in C/I.

# Block 1 - rows 3-8 (category labels unchanged, just add the new columns)
$ws.Range("E3").Value = "start DD"
$ws.Range("K3").Value = "start SD"

# Block 2 - rows 10-15
$ws.Range("C10").Value = "no walk/diff"
$ws.Range("E10").Value = "start SD"
$ws.Range("I10").Value = "no walk/diff"
$ws.Range("K10").Value = "start SD"

# Block 3 - rows 17-22
$ws.Range("C17").Value = "walk/diff"
$ws.Range("E17").Value = "start SD"
$ws.Range("I17").Value = "no walk/same"
$ws.Range("K17").Value = "start DD"

# Block 4 - rows 24-29
$ws.Range("C24").Value = "walk/same"
$ws.Range("E24").Value = "start DD"
$ws.Range("K24").Value = "start DD"

# Block 5 - rows 31-36
$ws.Range("C31").Value = "no walk/same"
$ws.Range("E31").Value = "start DD"
$ws.Range("K31").Value = "start SD"

# Cosmetic touches to mirror the saved file: last selected cell and page
# orientation.
$ws.PageSetup.Orientation = 1
$ws.Range("K33").Select()
